$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "=D(n-1)" running-carry formulas to column C, rows 210-225 ---
# (mirrors the two shared-formula groups C210:C222 / C223:C225 in the diff;
#  the underlying engine doesn't need the shared-formula optimisation, just
#  the equivalent per-cell formula + resulting value)
for ($r = 210; $r -le 225; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=D$prev"
}

# --- Bump column D values for rows 212-224 (these feed the C formulas above) ---
$ws.Range("D212").Value = 102.755
$ws.Range("D213").Value = 103.80500000000001
$ws.Range("D214").Value = 103.995
$ws.Range("D215").Value = 104.005
$ws.Range("D216").Value = 104.28
$ws.Range("D217").Value = 104.75
$ws.Range("D218").Value = 104.995
$ws.Range("D219").Value = 105.54
$ws.Range("D220").Value = 105.77500000000001
$ws.Range("D221").Value = 105.84
$ws.Range("D222").Value = 105.64
$ws.Range("D223").Value = 105.77
$ws.Range("D224").Value = 106.48

# --- Unrelated single-cell correction further down the sheet ---
$ws.Range("D230").Value = 112.28

# --- View/selection tweaks captured in the diff ---
$ws.Application.ActiveWindow.Zoom = 143
$ws.Range("E236").Select() | Out-Null
